$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" (column G) values per updated save_data (std/mean recalculated, s_vals rewritten)
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 3
$ws.Range("G8").Value = 1
